$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update U and V columns (previously 0) with new values, and AG totals accordingly.

$ws.Range("U2").Value = 21383.3
$ws.Range("V2").Value = 13491.76
$ws.Range("AG2").Value = 238595.79

$ws.Range("U3").Value = 8356
$ws.Range("V3").Value = 8011.9
$ws.Range("AG3").Value = 138572.6

$ws.Range("U4").Value = 2017
$ws.Range("V4").Value = 1656
$ws.Range("AG4").Value = 70572.60000000001

$ws.Range("U5").Value = 1334
$ws.Range("V5").Value = 1527
$ws.Range("AG5").Value = 65508.19

$ws.Range("U6").Value = 33090.3
$ws.Range("V6").Value = 24686.66
$ws.Range("AG6").Value = 513249.18
